# Rename header labels on the existing sheets
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "PO Forecast"

# Match the outline properties used on the other sheets (summaryBelow/summaryRight)
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Header row
$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Match header styling (bold, bordered, centered) used on the other sheets
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

# Match the date-formatted style used for column A on the other sheets
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A23").PasteSpecial(-4122)

# Data rows
$newSheet.Range("A2").Value = 45515.99999999999
$newSheet.Range("B2").Value = 52
$newSheet.Range("C2").Value = 9.989153854791459
$newSheet.Range("D2").Value = 96.03507057254475
$newSheet.Range("A3").Value = 45550.99999999999
$newSheet.Range("B3").Value = 49
$newSheet.Range("C3").Value = 8.95280658423493
$newSheet.Range("D3").Value = 93.67099430125106
$newSheet.Range("A4").Value = 45557.99999999999
$newSheet.Range("B4").Value = 49
$newSheet.Range("C4").Value = 9.021381237352241
$newSheet.Range("D4").Value = 90.254916181821
$newSheet.Range("A5").Value = 45564.99999999999
$newSheet.Range("B5").Value = 48
$newSheet.Range("C5").Value = 9.427265539656338
$newSheet.Range("D5").Value = 90.60135283742832
$newSheet.Range("A6").Value = 45571.99999999999
$newSheet.Range("B6").Value = 48
$newSheet.Range("C6").Value = 6.742444311963943
$newSheet.Range("D6").Value = 89.03375081904797
$newSheet.Range("A7").Value = 45578.99999999999
$newSheet.Range("B7").Value = 47
$newSheet.Range("C7").Value = 2.727989591982412
$newSheet.Range("D7").Value = 88.5129023331511
$newSheet.Range("A8").Value = 45585.99999999999
$newSheet.Range("B8").Value = 47
$newSheet.Range("C8").Value = 5.455404249151727
$newSheet.Range("D8").Value = 90.4242218632554
$newSheet.Range("A9").Value = 45592.99999999999
$newSheet.Range("B9").Value = 46
$newSheet.Range("C9").Value = 4.491031205301445
$newSheet.Range("D9").Value = 86.57789740745476
$newSheet.Range("A10").Value = 45599.99999999999
$newSheet.Range("B10").Value = 46
$newSheet.Range("C10").Value = 2.658746655745986
$newSheet.Range("D10").Value = 87.31032380839399
$newSheet.Range("A11").Value = 45606.99999999999
$newSheet.Range("B11").Value = 45
$newSheet.Range("C11").Value = -1.875962760761451
$newSheet.Range("D11").Value = 87.89307421002461
$newSheet.Range("A12").Value = 45613.99999999999
$newSheet.Range("B12").Value = 45
$newSheet.Range("C12").Value = 5.155799146327224
$newSheet.Range("D12").Value = 86.56155991716489
$newSheet.Range("A13").Value = 45627.99999999999
$newSheet.Range("B13").Value = 44
$newSheet.Range("C13").Value = 5.99059676563966
$newSheet.Range("D13").Value = 86.62274695521602
$newSheet.Range("A14").Value = 45634.99999999999
$newSheet.Range("B14").Value = 43
$newSheet.Range("C14").Value = 3.912920538875787
$newSheet.Range("D14").Value = 87.88064735762822
$newSheet.Range("A15").Value = 45641.99999999999
$newSheet.Range("B15").Value = 42
$newSheet.Range("C15").Value = -2.358907683166315
$newSheet.Range("D15").Value = 84.16554984552113
$newSheet.Range("A16").Value = 45648.99999999999
$newSheet.Range("B16").Value = 42
$newSheet.Range("C16").Value = -1.396820779576964
$newSheet.Range("D16").Value = 85.85998875022406
$newSheet.Range("A17").Value = 45655.99999999999
$newSheet.Range("B17").Value = 41
$newSheet.Range("C17").Value = -1.822809202396947
$newSheet.Range("D17").Value = 82.16440236740938
$newSheet.Range("A18").Value = 45662.99999999999
$newSheet.Range("B18").Value = 41
$newSheet.Range("C18").Value = 0.5042383975073585
$newSheet.Range("D18").Value = 81.55440441293351
$newSheet.Range("A19").Value = 45669.99999999999
$newSheet.Range("B19").Value = 40
$newSheet.Range("C19").Value = -3.830508451336271
$newSheet.Range("D19").Value = 80.01654293659945
$newSheet.Range("A20").Value = 45676.99999999999
$newSheet.Range("B20").Value = 40
$newSheet.Range("C20").Value = -1.456520232824235
$newSheet.Range("D20").Value = 82.81794883957481
$newSheet.Range("A21").Value = 45683.99999999999
$newSheet.Range("B21").Value = 39
$newSheet.Range("C21").Value = 0.16926601540088
$newSheet.Range("D21").Value = 80.31919734300132
$newSheet.Range("A22").Value = 45690.99999999999
$newSheet.Range("B22").Value = 39
$newSheet.Range("C22").Value = -2.702628567532695
$newSheet.Range("D22").Value = 79.13332760804872
$newSheet.Range("A23").Value = 45697.99999999999
$newSheet.Range("B23").Value = 38
$newSheet.Range("C23").Value = -3.619521353960328
$newSheet.Range("D23").Value = 80.11921535651275